$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

$ws.Cells.Item($row, 1).Value = "CUFR7F"
$ws.Cells.Item($row, 2).Value = "Rodillo de recogida de papel RM1-1497-000 para HP"
$ws.Cells.Item($row, 3).Value = "P1500 P1505 P1560 P1566 P1600 P1606 M1120 M1522 M1530 M1536 MF211 MF212 MF215 MF216 MF217 MF221 MF222 MF223 MF224 MF226 MF227 MF229 MF236 MF237 MF244 MF247 MF249"
$ws.Cells.Item($row, 4).Value = 6500
$ws.Cells.Item($row, 5).Value = 50000
$ws.Cells.Item($row, 6).Value = 9
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E83-D83)*G83"
$ws.Cells.Item($row, 9).Formula = "=D83*F83"
$ws.Cells.Item($row, 10).Value = 58500
